# Update "想去人数" (number of people interested) counts on the
# "展览" sheet and the "全部类型" sheet:
#   F2: 5412 -> 5423
#   F4: 930  -> 932

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5423
    $ws.Range("F4").Value = 932
}
